$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Active sheet: $($ws.Name)"
foreach ($sheet in $wb.Worksheets) {
    Write-Host $sheet.Name
}
